# Generate Report for Handoff
# - Flip status from "In Translation" to "Ready for handoff"
# - Refresh the handoff timestamps
# - Widen the "Status" columns so the new, longer status text fits

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# A ColumnWidth value that Excel's pixel-snapping rounds to the same stored
# column width (~17.216 "characters") used in the target report.
$statusColumnWidth = 16.33

# ---------------------------------------------------------------------
# Overview sheet: per-language status columns (E = zh-cn, F = de-de) and
# the overall "Latest HO Xliff Generate Date" column (G)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-11-23 11:58:12"

$wsOverview.Columns.Item(5).ColumnWidth = $statusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColumnWidth

# ---------------------------------------------------------------------
# zh-cn sheet: Status (C) and Latest Handoff Datetime (H)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-11-23 11:57:57"
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColumnWidth

# ---------------------------------------------------------------------
# de-de sheet: Status (C) and Latest Handoff Datetime (H)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-11-23 11:58:12"
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColumnWidth
